# Updated cryptos list - apply per-cell price/volume changes from upstream scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.723.22"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "1.847.29"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.80"
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4293"
$ws.Range("E7").Value = "  +1.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3654"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.05"
$ws.Range("E9").Value = "  -1.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07343"
$ws.Range("E10").Value = "  +0.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8760"
$ws.Range("E11").Value = "  -2.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.75"
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("D13").Value = "1.817.49"
$ws.Range("E13").Value = "  +0.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.343"
$ws.Range("E14").Value = "  -0.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.525"
$ws.Range("E15").Value = "  -0.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06941"
$ws.Range("E16").Value = "  +1.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "79.83"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000009018"
$ws.Range("E19").Value = "  +1.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.39"
$ws.Range("E21").Value = "  -1.14%  "
$ws.Range("D22").Value = "27.582.44"
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.969"
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.36"
$ws.Range("E24").Value = "  -2.36%  "
$ws.Range("D25").Value = "2.084.08"
$ws.Range("E25").Value = "  +4.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.978"
$ws.Range("E26").Value = "  -3.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.20"
$ws.Range("E27").Value = "  +1.22%  "
$ws.Range("E28").Value = "  +1.78%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.245"
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.15"
$ws.Range("E30").Value = "  +7.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.868"
$ws.Range("E31").Value = "  +1.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08897"
$ws.Range("E32").Value = "  +0.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7517"
$ws.Range("E33").Value = "  -3.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.548"
$ws.Range("E34").Value = "  -0.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.974"
$ws.Range("E35").Value = "  +1.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.123"
$ws.Range("E36").Value = "  +2.50%  "
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05425"
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.105"
$ws.Range("E39").Value = "  +0.98%  "
$ws.Range("E40").Value = "  +0.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.831"
$ws.Range("E41").Value = "  +0.49%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5085"
$ws.Range("E42").Value = "  +0.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1656"
$ws.Range("E43").Value = "  +0.82%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.593"
$ws.Range("E44").Value = "  -3.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.341"
$ws.Range("E45").Value = "  +1.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.06540"
$ws.Range("E46").Value = "  -1.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.37"
$ws.Range("E47").Value = "  +0.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "105.08"
$ws.Range("E48").Value = "  -0.79%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4659"
$ws.Range("E49").Value = "  -1.44%  "
$ws.Range("E50").Value = "  +0.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.625"
$ws.Range("E51").Value = "  -1.05%  "
